$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("Realizar login no aplicativo m?vel.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
Write-Host "found:" $found
Write-Host "before size:" $rng.Font.Size
$rng.Font.Size = 12
Write-Host "after size:" $rng.Font.Size
